$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 144.783305
$ws.Range("H2").Value = 434.349915
$ws.Range("I2").Value = 0.2430046335191003
$ws.Range("J2").Value = 0.251012682214973
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3372215
$ws.Range("N2").Value = 0.674443
$ws.Range("Q2").Value = 48.8240432870575
$ws.Range("R2").Value = 292.944259722345
$ws.Range("S2").Value = 0.2430046335191003
$ws.Range("T2").Value = 0.251012682214973

$ws.Range("G3").Value = 82.24887099999999
$ws.Range("I3").Value = 0.1380466950572427
$ws.Range("J3").Value = 0.1425959278859072
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3372215
$ws.Range("N3").Value = 0.674443
$ws.Range("Q3").Value = 27.7360876519265
$ws.Range("R3").Value = 166.416525911559
$ws.Range("S3").Value = 0.1380466950572427
$ws.Range("T3").Value = 0.1425959278859072

$ws.Range("G4").Value = 163.8590903333333
$ws.Range("H4").Value = 491.577271
$ws.Range("I4").Value = 0.2750214756820535
$ws.Range("J4").Value = 0.284084617144743
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3372215
$ws.Range("N4").Value = 0.674443
$ws.Range("Q4").Value = 55.25680823084217
$ws.Range("R4").Value = 331.540849385053
$ws.Range("S4").Value = 0.2750214756820535
$ws.Range("T4").Value = 0.284084617144743

$ws.Range("G5").Value = 57.0238095
$ws.Range("H5").Value = 114.047619
$ws.Range("I5").Value = 0.09570889357312636
$ws.Range("J5").Value = 0.06590860906562239
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3372215
$ws.Range("N5").Value = 0.674443
$ws.Range("Q5").Value = 19.22965457530425
$ws.Range("R5").Value = 76.91861830121699
$ws.Range("S5").Value = 0.09570889357312636
$ws.Range("T5").Value = 0.06590860906562239

$ws.Range("G6").Value = 147.8896333333333
$ws.Range("H6").Value = 443.6689
$ws.Range("I6").Value = 0.2482183021684772
$ws.Range("J6").Value = 0.2563981636887546
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3372215
$ws.Range("N6").Value = 0.674443
$ws.Range("Q6").Value = 49.87156398711667
$ws.Range("R6").Value = 299.2293839227
$ws.Range("S6").Value = 0.2482183021684772
$ws.Range("T6").Value = 0.2563981636887546
